$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 237
$ws1.Range("F3").Value = 1338
$ws1.Range("F5").Value = 886
$ws1.Range("F6").Value = 31
$ws1.Range("F7").Value = 1200
$ws1.Range("F8").Value = 1504
$ws1.Range("F9").Value = 151
$ws1.Range("F11").Value = 592
$ws1.Range("F12").Value = 429
$ws1.Range("F13").Value = 95
$ws1.Range("F16").Value = 84
$ws1.Range("F17").Value = 77
$ws1.Range("F18").Value = 5924
$ws1.Range("F20").Value = 5757
$ws1.Range("F21").Value = 9749
$ws1.Range("F22").Value = 120
$ws1.Range("F24").Value = 174
$ws1.Range("F25").Value = 264
$ws1.Range("F26").Value = 483
$ws1.Range("F27").Value = 160
$ws1.Range("F28").Value = 139
$ws1.Range("F29").Value = 4357
$ws1.Range("F30").Value = 354

# Sheet: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 609

# Sheet: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 609
$ws4.Range("F4").Value = 237
$ws4.Range("F5").Value = 1338
$ws4.Range("F8").Value = 887
$ws4.Range("F9").Value = 31
$ws4.Range("F10").Value = 1200
$ws4.Range("F12").Value = 1504
$ws4.Range("F14").Value = 151
$ws4.Range("F15").Value = 593
$ws4.Range("F17").Value = 429
$ws4.Range("F18").Value = 95
$ws4.Range("F22").Value = 84
$ws4.Range("F23").Value = 77
$ws4.Range("F24").Value = 5924
$ws4.Range("F26").Value = 5757
$ws4.Range("F27").Value = 9749
$ws4.Range("F29").Value = 120
$ws4.Range("F31").Value = 174
$ws4.Range("F32").Value = 264
$ws4.Range("F34").Value = 483
$ws4.Range("F37").Value = 160
$ws4.Range("F38").Value = 139
$ws4.Range("F39").Value = 4357
$ws4.Range("F46").Value = 354
